$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old long-form Portuguese "Programa" paragraph row.
# This single deletion realigns all subsequent rows/labels to match the target layout.
$ws.Rows.Item(16).Delete()

# Update cell contents that changed as part of the edit
$ws.Range('B10').Value = '6007846 - Júlio César dos Santos'
$ws.Range('C10').Value = '6007846 - Júlio César dos Santos'
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'
$ws.Range('C14').Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '01/01/2018'
$ws.Range('C15').Value = '01/01/2018'
$ws.Range('B18').Value = '6007846 - Júlio César dos Santos'
$ws.Range('C18').Value = '6007846 - Júlio César dos Santos'
$ws.Range('B19').Value = 'Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2).'
$ws.Range('C19').Value = 'Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2).'
$ws.Range('B20').Value = 'A média aritmética das notas corresponderá à média do período letivo, ou seja:Média do período letivo normal = (P1+ P2)/2.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0.'
$ws.Range('C20').Value = 'A média aritmética das notas corresponderá à média do período letivo, ou seja:Média do período letivo normal = (P1+ P2)/2.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0.'
$ws.Range('B21').Value = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'
$ws.Range('C21').Value = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'

# Fix up row heights that differ from what a plain row delete leaves behind
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120

Write-Host "Edit complete"
